$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the "Date" and "Jurisdiction" values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-09-15T20:59:49+00:00"
$meta.Range("B12").Value = "Global (Whole world)"

# --- Re-apply the wrap/vertical-top alignment on every used cell in both
#     sheets. The cells already render with top-aligned, wrapped text, but
#     doing this explicitly makes sure the alignment is flagged as applied
#     (applyAlignment) on the styles backing these cells, on every sheet -
#     not just the ones we touched above. ---
foreach ($sheetName in @("Metadata", "Concepts")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $used = $sheet.UsedRange
    $used.VerticalAlignment = -4160   # xlTop
    $used.WrapText = $true
}
